$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.030.13'
$ws.Range("E2").Value = '  +0.52%  '

# Row 3
$ws.Range("D3").Value = '1.683.29'
$ws.Range("E3").Value = '  +0.60%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.04%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.517'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.58%  '

# Row 7
$ws.Range("E7").Value = '  -0.06%  '

# Row 8
$ws.Range("B8").Value = 'Solana'
$ws.Range("C8").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.38'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.07%  '

# Row 9
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.252'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.42%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0621'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.02%  '

# Row 11
$ws.Range("E11").Value = '  -0.48%  '

# Row 12
$ws.Range("D12").Value = '1.919.77'
$ws.Range("E12").Value = '  +0.52%  '

# Row 13
$ws.Range("D13").Value = '1.680.96'
$ws.Range("E13").Value = '  +0.81%  '

# Row 14
$ws.Range("E14").Value = '  +0.39%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.535'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.14%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.14'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.61%  '

# Row 17
$ws.Range("D17").Value = '27.031.92'
$ws.Range("E17").Value = '  +0.28%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.17'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.59%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '236.55'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.45%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0736'
$ws.Range("E20").Value = '  +0.02%  '

# Row 21
$ws.Range("E21").Value = '  -0.02%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.47'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.07%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.41%  '

# Row 24
$ws.Range("E24").Value = '  -4.00%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.85'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.80%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.17%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.70%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.113'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.90%  '

# Row 29
$ws.Range("E29").Value = '  -0.02%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0501'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.63%  '

# Row 31
$ws.Range("E31").Value = '  -0.34%  '

# Row 32
$ws.Range("E32").Value = '  +0.10%  '

# Row 33
$ws.Range("D33").Value = '1.515.68'
$ws.Range("E33").Value = '  +3.64%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.19'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.57%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.69'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.39%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.40'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.51%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.590'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.42%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.922'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.49%  '

# Row 39
$ws.Range("E39").Value = '  +3.00%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.05'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.56%  '

# Row 41
$ws.Range("E41").Value = '  -5.38%  '

# Row 42
$ws.Range("E42").Value = '  +0.04%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '68.33'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.77%  '

# Row 44
$ws.Range("E44").Value = '  -1.16%  '

# Row 45
$ws.Range("D45").Value = '1.825.09'
$ws.Range("E45").Value = '  +0.12%  '

# Row 46
$ws.Range("E46").Value = '  +0.29%  '

# Row 47
$ws.Range("E47").Value = '  -0.42%  '

# Row 48
$ws.Range("E48").Value = '  +4.18%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.53'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.72%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.87'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.20%  '

# Row 51
$ws.Range("E51").Value = '  +0.16%  '
